$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'301.62"
$ws.Range("E2").Value = "'1.24%"
$ws.Range("G2").Value = "'14"
$ws.Range("D3").Value = "'31.93"
$ws.Range("E3").Value = "'2.20%"
$ws.Range("G3").Value = "'14"
$ws.Range("D4").Value = "'5.121"
$ws.Range("E4").Value = "'0.39%"
$ws.Range("G4").Value = "'14"
$ws.Range("D5").Value = "'0.07848"
$ws.Range("E5").Value = "'-1.95%"
$ws.Range("G5").Value = "'14"
$ws.Range("D6").Value = "'2.263"
$ws.Range("E6").Value = "'-8.26%"
$ws.Range("G6").Value = "'14"
$ws.Range("D7").Value = "'7.809"
$ws.Range("E7").Value = "'-0.01%"
$ws.Range("G7").Value = "'14"
$ws.Range("E8").Value = "'0.08%"
$ws.Range("G8").Value = "'14"
$ws.Range("E9").Value = "'0.85%"
$ws.Range("G9").Value = "'14"
$ws.Range("D10").Value = "'0.1774"
$ws.Range("E10").Value = "'2.39%"
$ws.Range("G10").Value = "'14"
$ws.Range("D11").Value = "'0.07691"
$ws.Range("E11").Value = "'5.46%"
$ws.Range("G11").Value = "'14"
$ws.Range("D12").Value = "'0.08873"
$ws.Range("E12").Value = "'4.13%"
$ws.Range("G12").Value = "'14"
$ws.Range("D13").Value = "'0.03098"
$ws.Range("E13").Value = "'1.79%"
$ws.Range("G13").Value = "'14"
$ws.Range("E14").Value = "'0.50%"
$ws.Range("G14").Value = "'14"
$ws.Range("D15").Value = "'0.001510"
$ws.Range("E15").Value = "'1.15%"
$ws.Range("G15").Value = "'14"
$ws.Range("D16").Value = "'0.006018"
$ws.Range("E16").Value = "'-0.60%"
$ws.Range("G16").Value = "'14"
$ws.Range("D17").Value = "'3.468"
$ws.Range("E17").Value = "'-1.34%"
$ws.Range("G17").Value = "'14"
$ws.Range("E18").Value = "'0.09%"
$ws.Range("G18").Value = "'14"
$ws.Range("E19").Value = "'0.22%"
$ws.Range("G19").Value = "'14"
$ws.Range("E20").Value = "'-1.40%"
$ws.Range("G20").Value = "'14"
$ws.Range("D21").Value = "'4.331"
$ws.Range("E21").Value = "'-6.22%"
$ws.Range("G21").Value = "'14"
$ws.Range("E22").Value = "'10.53%"
$ws.Range("G22").Value = "'14"
$ws.Range("D23").Value = "'0.04598"
$ws.Range("E23").Value = "'-1.00%"
$ws.Range("G23").Value = "'14"
$ws.Range("D24").Value = "'0.001252"
$ws.Range("E24").Value = "'0.27%"
$ws.Range("G24").Value = "'14"
$ws.Range("E25").Value = "'0.99%"
$ws.Range("G25").Value = "'14"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'4.01%"
$ws.Range("G26").Value = "'14"
$ws.Range("E27").Value = "'-1.40%"
$ws.Range("G27").Value = "'14"
$ws.Range("G28").Value = "'14"
$ws.Range("G29").Value = "'14"
$ws.Range("G30").Value = "'14"
$ws.Range("G31").Value = "'14"
$ws.Range("G32").Value = "'14"
$ws.Range("G33").Value = "'14"
$ws.Range("G34").Value = "'14"
$ws.Range("G35").Value = "'14"
$ws.Range("G36").Value = "'14"
$ws.Range("G37").Value = "'14"
$ws.Range("G38").Value = "'14"
$ws.Range("D39").Value = "'0.01787"
$ws.Range("E39").Value = "'0.14%"
$ws.Range("G39").Value = "'14"
$ws.Range("D40").Value = "'0.04802"
$ws.Range("E40").Value = "'7.84%"
$ws.Range("G40").Value = "'14"
$ws.Range("D41").Value = "'0.007253"
$ws.Range("E41").Value = "'4.18%"
$ws.Range("G41").Value = "'14"
$ws.Range("D42").Value = "'0.1365"
$ws.Range("E42").Value = "'1.75%"
$ws.Range("G42").Value = "'14"
$ws.Range("D43").Value = "'0.002189"
$ws.Range("E43").Value = "'-2.37%"
$ws.Range("G43").Value = "'14"
$ws.Range("D44").Value = "'0.01134"
$ws.Range("E44").Value = "'15.75%"
$ws.Range("G44").Value = "'14"
$ws.Range("D45").Value = "'0.00006245"
$ws.Range("E45").Value = "'-5.38%"
$ws.Range("G45").Value = "'14"
$ws.Range("E46").Value = "'-0.13%"
$ws.Range("G46").Value = "'14"
$ws.Range("G47").Value = "'14"
$ws.Range("D48").Value = "'0.6972"
$ws.Range("E48").Value = "'-15.03%"
$ws.Range("G48").Value = "'14"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("G49").Value = "'14"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("G50").Value = "'14"
$ws.Range("G51").Value = "'14"
